$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.707.55"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.077.02"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'233.83"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'58.17"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "'0.391"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'0.0782"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.383.75"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'14.86"
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "'20.90"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "'0.772"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "2.057.61"
$ws.Range("E17").Value = "  -2.91%  "
$ws.Range("D18").Value = "37.693.49"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'6.17"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "'71.19"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").Value = "'227.90"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "'169.55"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("E27").Value = "  +3.96%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'1.40"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("E31").Value = "  +2.19%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "'4.67"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "'5.35"
$ws.Range("E39").Value = "  -4.99%  "
$ws.Range("D40").Value = "'0.0976"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").Value = "'98.07"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "1.453.92"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").Value = "'16.57"
$ws.Range("E45").Value = "  +6.42%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").Value = "'4.24"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").Value = "'7.39"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "2.267.66"
$ws.Range("E51").Value = "  -1.59%  "
